$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 600
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 600
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -826
$ws.Range("H4").Value = 339.7143
$ws.Range("I4").Value = 229.66667
$ws.Range("K4").Value = 229.66667
$ws.Range("M4").Value = -115.66667
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H9").Value = 165
$ws.Range("I9").Value = 29
$ws.Range("J9").Value = 210.33333
$ws.Range("K9").Value = 29
$ws.Range("L9").Value = 210.33333
$ws.Range("M9").Value = 140
$ws.Range("N9").Value = -548.3333299999999
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H28").Value = 42648.625
$ws.Range("I28").Value = 50558.6
$ws.Range("J28").Value = 3098.75
$ws.Range("K28").Value = 50558.6
$ws.Range("L28").Value = 3098.75
$ws.Range("M28").Value = -50073.6
$ws.Range("N28").Value = -4068.75
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H93").Value = 50000
$ws.Range("J93").Value = 50000
$ws.Range("L93").Value = 50000
$ws.Range("N93").Value = -54992
$ws.Range("H132").Value = 1350.0984
$ws.Range("I132").Value = 1334.2
$ws.Range("K132").Value = 4002.6
$ws.Range("M132").Value = -1472.6
$ws.Range("H138").Value = 3411.7632
$ws.Range("I138").Value = 4665.1665
$ws.Range("J138").Value = 3176.75
$ws.Range("K138").Value = 13995.4995
$ws.Range("L138").Value = 9530.25
$ws.Range("M138").Value = -8855.499500000002
$ws.Range("N138").Value = -19810.25
$ws.Range("H141").Value = 1095
$ws.Range("I141").Value = 1095
$ws.Range("K141").Value = 3285
$ws.Range("M141").Value = 1895
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6338.1
$ws.Range("I2").Value = 831.9091
$ws.Range("K2").Value = 831.9091
$ws.Range("M2").Value = -718.9091
$ws.Range("H32").Value = 2707
$ws.Range("I32").Value = 2707
$ws.Range("K32").Value = 2707
$ws.Range("M32").Value = -2420
$ws.Range("H53").Value = 46292.668
$ws.Range("I53").Value = 38888
$ws.Range("K53").Value = 38888
$ws.Range("M53").Value = -38206
$ws.Range("H116").Value = 6338.1
$ws.Range("I116").Value = 831.9091
$ws.Range("K116").Value = 831.9091
$ws.Range("M116").Value = 1462.0909
$ws.Range("H119").Value = 66268.5
$ws.Range("J119").Value = 66268.5
$ws.Range("L119").Value = 66268.5
$ws.Range("N119").Value = -75944.5
$ws.Range("H132").Value = 8379.097
$ws.Range("I132").Value = 4601.9
$ws.Range("K132").Value = 13805.7
$ws.Range("M132").Value = -11275.7
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6338.1
$ws.Range("I3").Value = 831.9091
$ws.Range("K3").Value = 831.9091
$ws.Range("M3").Value = -717.9091
$ws.Range("H12").Value = 4410
$ws.Range("I12").Value = 150
$ws.Range("J12").Value = 5475
$ws.Range("K12").Value = 150
$ws.Range("L12").Value = 5475
$ws.Range("M12").Value = 18
$ws.Range("N12").Value = -5811
$ws.Range("H20").Value = 3736.1155
$ws.Range("I20").Value = 2871.4
$ws.Range("J20").Value = 4915.273
$ws.Range("K20").Value = 2871.4
$ws.Range("L20").Value = 4915.273
$ws.Range("M20").Value = -2624.4
$ws.Range("N20").Value = -5409.273
$ws.Range("H64").Value = 990
$ws.Range("I64").Value = 1000
$ws.Range("J64").Value = 980
$ws.Range("K64").Value = 1000
$ws.Range("L64").Value = 980
$ws.Range("M64").Value = -775
$ws.Range("N64").Value = -1430
$ws.Range("H67").Value = 990
$ws.Range("I67").Value = 1000
$ws.Range("J67").Value = 980
$ws.Range("K67").Value = 1000
$ws.Range("L67").Value = 980
$ws.Range("M67").Value = -220
$ws.Range("N67").Value = -2540
$ws.Range("H86").Value = 3107
$ws.Range("I86").Value = 3107
$ws.Range("K86").Value = 3107
$ws.Range("M86").Value = -1984
$ws.Range("H89").Value = 3107
$ws.Range("I89").Value = 3107
$ws.Range("K89").Value = 15535
$ws.Range("M89").Value = -9919
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 20000
$ws.Range("J45").Value = 20000
$ws.Range("L45").Value = 20000
$ws.Range("N45").Value = -21186
$ws.Range("H58").Value = 5101.5264
$ws.Range("I58").Value = 2155.0908
$ws.Range("K58").Value = 2155.0908
$ws.Range("M58").Value = -1952.0908
$ws.Range("H86").Value = 11847.5
$ws.Range("J86").Value = 15833.667
$ws.Range("L86").Value = 15833.667
$ws.Range("N86").Value = -18079.667
$ws.Range("H89").Value = 11847.5
$ws.Range("J89").Value = 15833.667
$ws.Range("L89").Value = 79168.33499999999
$ws.Range("N89").Value = -90400.33499999999
$ws.Range("H136").Value = 5101.5264
$ws.Range("I136").Value = 2155.0908
$ws.Range("K136").Value = 6465.2724
$ws.Range("M136").Value = -3915.2724
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 5187810
$ws.Range("J34").Value = 8334.333000000001
$ws.Range("L34").Value = 25002.999
$ws.Range("N34").Value = -25170.999
$ws.Range("H39").Value = 3325
$ws.Range("J39").Value = 2900
$ws.Range("L39").Value = 8700
$ws.Range("N39").Value = -9288
$ws.Range("H55").Value = 2118.5
$ws.Range("J55").Value = 5800
$ws.Range("L55").Value = 17400
$ws.Range("N55").Value = -17754
$ws.Range("H87").Value = 16499.5
$ws.Range("I87").Value = 16499.5
$ws.Range("K87").Value = 49498.5
$ws.Range("M87").Value = -48250.5
$ws.Range("H90").Value = 16499.5
$ws.Range("I90").Value = 16499.5
$ws.Range("K90").Value = 148495.5
$ws.Range("M90").Value = -142255.5
$ws.Range("H140").Value = 4874.25
$ws.Range("I140").Value = 1750
$ws.Range("K140").Value = 5250
$ws.Range("M140").Value = -70
$ws.Range("H141").Value = 6921.357
$ws.Range("I141").Value = 3618.182
$ws.Range("K141").Value = 10854.546
$ws.Range("M141").Value = -5674.545999999998
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11466.429
$ws.Range("I70").Value = 8351.75
$ws.Range("J70").Value = 15619.333
$ws.Range("K70").Value = 8351.75
$ws.Range("L70").Value = 15619.333
$ws.Range("M70").Value = -8081.75
$ws.Range("N70").Value = -16159.333
$ws.Range("H73").Value = 11466.429
$ws.Range("I73").Value = 8351.75
$ws.Range("J73").Value = 15619.333
$ws.Range("K73").Value = 8351.75
$ws.Range("L73").Value = 15619.333
$ws.Range("M73").Value = -7415.75
$ws.Range("N73").Value = -17491.333
$ws.Range("H80").Value = 11088.25
$ws.Range("I80").Value = 7940
$ws.Range("J80").Value = 16335.333
$ws.Range("K80").Value = 7940
$ws.Range("L80").Value = 16335.333
$ws.Range("M80").Value = -6942
$ws.Range("N80").Value = -18331.333
$ws.Range("H83").Value = 11088.25
$ws.Range("I83").Value = 7940
$ws.Range("J83").Value = 16335.333
$ws.Range("K83").Value = 39700
$ws.Range("L83").Value = 81676.66500000001
$ws.Range("M83").Value = -34708
$ws.Range("N83").Value = -91660.66500000001
$ws.Range("H102").Value = 2873.5557
$ws.Range("I102").Value = 2454.353
$ws.Range("K102").Value = 2454.353
$ws.Range("M102").Value = -832.3530000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6792.4517
$ws.Range("I7").Value = 4058.48
$ws.Range("K7").Value = 4058.48
$ws.Range("M7").Value = -3946.48
$ws.Range("H40").Value = 7922.4287
$ws.Range("I40").Value = 9001
$ws.Range("K40").Value = 9001
$ws.Range("M40").Value = -8865
$ws.Range("H63").Value = 45000
$ws.Range("I63").Value = 45000
$ws.Range("K63").Value = 45000
$ws.Range("M63").Value = -44251
$ws.Range("H66").Value = 45000
$ws.Range("I66").Value = 45000
$ws.Range("K66").Value = 135000
$ws.Range("M66").Value = -131256
$ws.Range("H126").Value = 6792.4517
$ws.Range("I126").Value = 4058.48
$ws.Range("K126").Value = 12175.44
$ws.Range("M126").Value = -9705.440000000001
$ws.Range("H133").Value = 55326
$ws.Range("J133").Value = 55326
$ws.Range("L133").Value = 55326
$ws.Range("N133").Value = -60386
$ws.Range("H136").Value = 10067.81
$ws.Range("I136").Value = 4492.727
$ws.Range("J136").Value = 16200.4
$ws.Range("K136").Value = 13478.181
$ws.Range("L136").Value = 48601.2
$ws.Range("M136").Value = -10928.181
$ws.Range("N136").Value = -53701.2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 17615.666
$ws.Range("J41").Value = 18138.8
$ws.Range("L41").Value = 18138.8
$ws.Range("N41").Value = -18918.8
$ws.Range("H81").Value = 2976
$ws.Range("J81").Value = 9694.5
$ws.Range("L81").Value = 19389
$ws.Range("N81").Value = -21511
$ws.Range("H84").Value = 2976
$ws.Range("J84").Value = 9694.5
$ws.Range("L84").Value = 96945
$ws.Range("N84").Value = -107553
$ws.Range("H122").Value = 4630.9443
$ws.Range("I122").Value = 1853
$ws.Range("K122").Value = 5559
$ws.Range("M122").Value = -3109
